# Edit script for Extension-CareConnect-EmergencyCareDischargeStatus-1.xlsx
# Splits the single "Extension.valueCodeableConcept" row into:
#   Row 6: generic "Extension.value[x]" slicing-parent row (hidden)
#   Row 7: specific "valueCodeableConcept" slice row (new, hidden)
# and adjusts dimension / autofilter / conditional formatting / defined name / column width accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Update existing row 6 so it becomes the generic "Extension.value[x]" row
# ---------------------------------------------------------------------------
$ws.Range("A6").Value = "Extension.value[x]"
$ws.Range("B6").Value = ""
$ws.Range("E6").Value = "0"
$ws.Range("K6").Value = "Value of extension"
$ws.Range("L6").Value = "Value of extension - may be a resource or one of a constrained set of the data types (see Extensibility in the spec for list)."
$ws.Range("W6").Value = ""
$ws.Range("X6").Value = ""
$ws.Range("Y6").Value = ""
$ws.Range("AA6").Value = "type:`$this}`n"
$ws.Range("AB6").Value = ""
$ws.Range("AD6").Value = "closed"
# A6..AJ6 already correct for: C,D,F,G,H,I,J,M,N,O,P,Q,R,S,T,U,V,Z,AC,AE,AF,AG,AH,AI,AJ

# Undo any auto row-height recalculation triggered by the embedded newline above
$ws.Rows.Item(6).AutoFit()

# ---------------------------------------------------------------------------
# 2. Insert new row 7 (copy row 6 formatting) for the "valueCodeableConcept" slice
# ---------------------------------------------------------------------------
$ws.Rows.Item(7).Insert()
$ws.Range("A6:AJ6").Copy()
$ws.Range("A7:AJ7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A7").Value = "Extension.value[x]"
$ws.Range("B7").Value = "valueCodeableConcept"
$ws.Range("C7").Value = ""
$ws.Range("E7").Value = "1"
$ws.Range("F7").Value = "1"
$ws.Range("G7").Value = ""
$ws.Range("H7").Value = ""
$ws.Range("I7").Value = ""
$ws.Range("J7").Value = "CodeableConcept`n"
$ws.Range("K7").Value = "The status of the Patient on discharge from an Emergency Care Department."
$ws.Range("L7").Value = "Value of extension - may be a resource or one of a constrained set of the data types (see Extensibility in the spec for list)."
$ws.Range("O7").Value = ""
$ws.Range("Q7").Value = ""
$ws.Range("R7").Value = ""
$ws.Range("S7").Value = ""
$ws.Range("T7").Value = ""
$ws.Range("U7").Value = ""
$ws.Range("V7").Value = ""
$ws.Range("W7").Value = "required"
$ws.Range("X7").Value = "The status of the Patient on discharge from an Emergency Care Department."
$ws.Range("Y7").Value = "https://fhir.hl7.org.uk/STU3/ValueSet/CareConnect-EmergencyCareDischargeStatus-1"
$ws.Range("Z7").Value = ""
$ws.Range("AA7").Value = ""
$ws.Range("AB7").Value = ""
$ws.Range("AC7").Value = ""
$ws.Range("AD7").Value = ""
$ws.Range("AE7").Value = "Extension.value[x]"
$ws.Range("AF7").Value = "0"
$ws.Range("AG7").Value = "1"
$ws.Range("AH7").Value = ""
$ws.Range("AI7").Value = ""
$ws.Range("AJ7").Value = "N/A"

# Undo any auto row-height recalculation triggered by the embedded newline above
$ws.Rows.Item(7).AutoFit()

# ---------------------------------------------------------------------------
# 3. Re-assert the "hidden" state of the detail rows (2-7)
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 7; $r++) {
    $ws.Rows.Item($r).Hidden = $true
}

# ---------------------------------------------------------------------------
# 4. Narrow column A (30.890625 -> ~19.00390625 width units)
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 18.17

# ---------------------------------------------------------------------------
# 5. Fix up the autofilter range + re-apply the original filter criteria
# ---------------------------------------------------------------------------
$ws.AutoFilterMode = $false
$ws.Range("A1:AJ7").AutoFilter(7, "<>" + " ")
$blanks = @("")
$ws.AutoFilter.Range.AutoFilter(27, $blanks, 7)

# ---------------------------------------------------------------------------
# 6. Extend the conditional formatting range to include the new row
# ---------------------------------------------------------------------------
$fcs = $ws.Range("A2:AI5").FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($ws.Range("A2:AI6"))
}

# ---------------------------------------------------------------------------
# 7. Fix the defined name (_xlnm._FilterDatabase) range
# ---------------------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "=Elements!`$A`$1:`$AJ`$7"
    }
}
